# Add two more test records to the GBIF interpretation test dataset:
#   - a row validating depthAccuracy (inserted as the new row 61, right
#     after the existing "non numeric depth" depth-block rows)
#   - a row validating elevationAccuracy (inserted as the new row 68,
#     right after the existing elevation-block rows)
# Every row below each insertion point shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the depthAccuracy test row at row 61 -----------------------
$ws.Rows.Item(61).Insert() | Out-Null

$ws.Range("A61").Formula = '=CONCATENATE("urn:lsid:gbif.org:Test:",ROW(A61))'
$ws.Range("B61").Value = "Validate depthAccuracy calculated correctly from (max - min) / 2. Expected: depthAccuracy = 1"
$ws.Range("U61").Value = 12
$ws.Range("V61").Value = 10
$ws.Range("Y61").Value = "Puma concolor (Linnaeus, 1771)"

# --- Insert the elevationAccuracy test row at row 68 --------------------
# (old row 67 "taxon match FUZZY" is now at row 68, having been pushed
# down by the insert above; inserting here pushes it to row 69)
$ws.Rows.Item(68).Insert() | Out-Null

$ws.Range("A68").Formula = '=CONCATENATE("urn:lsid:gbif.org:Test:",ROW(A68))'
$ws.Range("B68").Value = "Validate elevationAccuracy calculated correctly from (max - min) / 2. Expected: elevationAccuracy = 1"
$ws.Range("W68").Value = 12
$ws.Range("X68").Value = 10
$ws.Range("Y68").Value = "Puma concolor (Linnaeus, 1771)"

# --- Match the saved selection state of the edited workbook -------------
$ws.Range("B68").Select() | Out-Null
